$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the column headers in row 1 so they carry the respective
#    input-file-name suffix instead of the generic "_old" / "_new" suffix.
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = "$($baseNames[$i])_FV2304"
}
# column 11 ("diff") keeps its name

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = "$($baseNames[$i])_FV2310"
}

# 2) Turn the used range into an Excel Table (ListObject) so the header row
#    gets the filter buttons / structured reference support.
$usedRange = $ws.UsedRange
$tbl = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row (split below row 1, top-left cell of the
#    scrollable area is A2).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
